$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed row 97 from row 96: first the cell formats (keeps the exact same
# style indices used throughout the sheet - bold/bordered for col A,
# the datetime number format for col E, default elsewhere), then the
# cell values/types (this also correctly seeds pais/torneio/temporada
# columns B/C/D as text, since those three columns are identical between
# row 96 and the new row 97).
$ws.Range("A96:V96").Copy()
$ws.Range("A97:V97").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A96:V96").Copy()
$ws.Range("A97:V97").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Now overwrite only the cells whose content actually differs for the
# new match (everything except B/C/D, which stay "ecuador"/"liga-pro"/"2023").
$ws.Range("A97").Value = 96
$ws.Range("E97").Value = 45239.04166666666
$ws.Range("F97").Value = "Tecnico U."
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = "LDU Quito"
$ws.Range("I97").Value = 2
$ws.Range("J97").Value = 2.47
$ws.Range("K97").Value = "05/11/2023 18:43"
$ws.Range("L97").Value = 3.14
$ws.Range("M97").Value = "09/11/2023 00:58"
$ws.Range("N97").Value = 3.37
$ws.Range("O97").Value = "05/11/2023 18:43"
$ws.Range("P97").Value = 3.29
$ws.Range("Q97").Value = "09/11/2023 00:56"
$ws.Range("R97").Value = 2.71
$ws.Range("S97").Value = "05/11/2023 18:43"
$ws.Range("T97").Value = 2.37
$ws.Range("U97").Value = "09/11/2023 00:58"
$ws.Range("V97").Value = "https://www.betexplorer.com/football/ecuador/liga-pro/tecnico-u-ldu-quito/8l0b4mjO/"
